# Daily Report update: 2026-01-15
# Appends the 2026-01-14 (Excel serial 46036) daily figures for each
# depository/Region_Type to Daily_Data, then refreshes the dependent
# Today_Summary and Monthly_Stats rollups for the depositories whose
# totals moved (BRINK'S, INC.; DELAWARE DEPOSITORY; MANFRA, TORDELLA &
# BROOKES, LLC).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Daily_Data: append rows 178-199 for date 46036 (2026-01-14)
# ---------------------------------------------------------------------
$wsDaily = $wb.Worksheets.Item("Daily_Data")

$newRowsData = New-Object 'object[,]' 22,8
$newRowsData[0,0] = 46036.0
$newRowsData[0,1] = "ASAHI DEPOSITORY LLC Registered"
$newRowsData[0,2] = 0.0
$newRowsData[0,3] = 0.0
$newRowsData[0,4] = 0.0
$newRowsData[0,5] = 0.0
$newRowsData[0,6] = 0.0
$newRowsData[0,7] = 0.0
$newRowsData[1,0] = 46036.0
$newRowsData[1,1] = "ASAHI DEPOSITORY LLC Eligible"
$newRowsData[1,2] = 0.0
$newRowsData[1,3] = 0.0
$newRowsData[1,4] = 0.0
$newRowsData[1,5] = 0.0
$newRowsData[1,6] = 0.0
$newRowsData[1,7] = 0.0
$newRowsData[2,0] = 46036.0
$newRowsData[2,1] = "BRINK'S, INC. Registered"
$newRowsData[2,2] = 90027.725
$newRowsData[2,3] = 0.0
$newRowsData[2,4] = 0.0
$newRowsData[2,5] = 0.0
$newRowsData[2,6] = 5490.05
$newRowsData[2,7] = 95517.775
$newRowsData[3,0] = 46036.0
$newRowsData[3,1] = "BRINK'S, INC. Eligible"
$newRowsData[3,2] = 5075.067
$newRowsData[3,3] = 0.0
$newRowsData[3,4] = 0.0
$newRowsData[3,5] = 0.0
$newRowsData[3,6] = 18635.207
$newRowsData[3,7] = 23710.274
$newRowsData[4,0] = 46036.0
$newRowsData[4,1] = "CNT DEPOSITORY, INC. Registered"
$newRowsData[4,2] = 1246.06
$newRowsData[4,3] = 0.0
$newRowsData[4,4] = 0.0
$newRowsData[4,5] = 0.0
$newRowsData[4,6] = 0.0
$newRowsData[4,7] = 1246.06
$newRowsData[5,0] = 46036.0
$newRowsData[5,1] = "CNT DEPOSITORY, INC. Eligible"
$newRowsData[5,2] = 0.0
$newRowsData[5,3] = 0.0
$newRowsData[5,4] = 0.0
$newRowsData[5,5] = 0.0
$newRowsData[5,6] = 0.0
$newRowsData[5,7] = 0.0
$newRowsData[6,0] = 46036.0
$newRowsData[6,1] = "DELAWARE DEPOSITORY Registered"
$newRowsData[6,2] = 1633.941
$newRowsData[6,3] = 0.0
$newRowsData[6,4] = 0.0
$newRowsData[6,5] = 0.0
$newRowsData[6,6] = 0.0
$newRowsData[6,7] = 1633.941
$newRowsData[7,0] = 46036.0
$newRowsData[7,1] = "DELAWARE DEPOSITORY Eligible"
$newRowsData[7,2] = 18509.729
$newRowsData[7,3] = 0.0
$newRowsData[7,4] = 50.145
$newRowsData[7,5] = -50.145
$newRowsData[7,6] = 0.0
$newRowsData[7,7] = 18459.584
$newRowsData[8,0] = 46036.0
$newRowsData[8,1] = "HSBC BANK, USA Registered"
$newRowsData[8,2] = 1295.223
$newRowsData[8,3] = 0.0
$newRowsData[8,4] = 0.0
$newRowsData[8,5] = 0.0
$newRowsData[8,6] = 0.0
$newRowsData[8,7] = 1295.223
$newRowsData[9,0] = 46036.0
$newRowsData[9,1] = "HSBC BANK, USA Eligible"
$newRowsData[9,2] = 9281.979
$newRowsData[9,3] = 0.0
$newRowsData[9,4] = 0.0
$newRowsData[9,5] = 0.0
$newRowsData[9,6] = 0.0
$newRowsData[9,7] = 9281.979
$newRowsData[10,0] = 46036.0
$newRowsData[10,1] = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered"
$newRowsData[10,2] = 2395.448
$newRowsData[10,3] = 0.0
$newRowsData[10,4] = 0.0
$newRowsData[10,5] = 0.0
$newRowsData[10,6] = 0.0
$newRowsData[10,7] = 2395.448
$newRowsData[11,0] = 46036.0
$newRowsData[11,1] = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible"
$newRowsData[11,2] = 0.0
$newRowsData[11,3] = 0.0
$newRowsData[11,4] = 0.0
$newRowsData[11,5] = 0.0
$newRowsData[11,6] = 0.0
$newRowsData[11,7] = 0.0
$newRowsData[12,0] = 46036.0
$newRowsData[12,1] = "JP MORGAN CHASE BANK NA Registered"
$newRowsData[12,2] = 124991.729
$newRowsData[12,3] = 0.0
$newRowsData[12,4] = 0.0
$newRowsData[12,5] = 0.0
$newRowsData[12,6] = 0.0
$newRowsData[12,7] = 124991.729
$newRowsData[13,0] = 46036.0
$newRowsData[13,1] = "JP MORGAN CHASE BANK NA Eligible"
$newRowsData[13,2] = 125407.673
$newRowsData[13,3] = 0.0
$newRowsData[13,4] = 0.0
$newRowsData[13,5] = 0.0
$newRowsData[13,6] = 0.0
$newRowsData[13,7] = 125407.673
$newRowsData[14,0] = 46036.0
$newRowsData[14,1] = "LOOMIS INTERNATIONAL (US) LLC Registered"
$newRowsData[14,2] = 68084.33
$newRowsData[14,3] = 0.0
$newRowsData[14,4] = 0.0
$newRowsData[14,5] = 0.0
$newRowsData[14,6] = 0.0
$newRowsData[14,7] = 68084.33
$newRowsData[15,0] = 46036.0
$newRowsData[15,1] = "LOOMIS INTERNATIONAL (US) LLC Eligible"
$newRowsData[15,2] = 116365.524
$newRowsData[15,3] = 0.0
$newRowsData[15,4] = 0.0
$newRowsData[15,5] = 0.0
$newRowsData[15,6] = 0.0
$newRowsData[15,7] = 116365.524
$newRowsData[16,0] = 46036.0
$newRowsData[16,1] = "MALCA-AMIT USA, LLC Registered"
$newRowsData[16,2] = 395.145
$newRowsData[16,3] = 0.0
$newRowsData[16,4] = 0.0
$newRowsData[16,5] = 0.0
$newRowsData[16,6] = 0.0
$newRowsData[16,7] = 395.145
$newRowsData[17,0] = 46036.0
$newRowsData[17,1] = "MALCA-AMIT USA, LLC Eligible"
$newRowsData[17,2] = 0.0
$newRowsData[17,3] = 0.0
$newRowsData[17,4] = 0.0
$newRowsData[17,5] = 0.0
$newRowsData[17,6] = 0.0
$newRowsData[17,7] = 0.0
$newRowsData[18,0] = 46036.0
$newRowsData[18,1] = "MANFRA, TORDELLA & BROOKES, LLC Registered"
$newRowsData[18,2] = 54605.27
$newRowsData[18,3] = 0.0
$newRowsData[18,4] = 0.0
$newRowsData[18,5] = 0.0
$newRowsData[18,6] = 5695.979
$newRowsData[18,7] = 60301.249
$newRowsData[19,0] = 46036.0
$newRowsData[19,1] = "MANFRA, TORDELLA & BROOKES, LLC Eligible"
$newRowsData[19,2] = 1068.408
$newRowsData[19,3] = 0.0
$newRowsData[19,4] = 0.0
$newRowsData[19,5] = 0.0
$newRowsData[19,6] = 0.0
$newRowsData[19,7] = 1068.408
$newRowsData[20,0] = 46036.0
$newRowsData[20,1] = "STONEX PRECIOUS METALS LLC Registered"
$newRowsData[20,2] = 14122.765
$newRowsData[20,3] = 0.0
$newRowsData[20,4] = 0.0
$newRowsData[20,5] = 0.0
$newRowsData[20,6] = 0.0
$newRowsData[20,7] = 14122.765
$newRowsData[21,0] = 46036.0
$newRowsData[21,1] = "STONEX PRECIOUS METALS LLC Eligible"
$newRowsData[21,2] = 16.075
$newRowsData[21,3] = 0.0
$newRowsData[21,4] = 0.0
$newRowsData[21,5] = 0.0
$newRowsData[21,6] = 0.0
$newRowsData[21,7] = 16.075

$wsDaily.Range("A178:H199").Value = $newRowsData
$wsDaily.Range("A178:A199").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# 2) Today_Summary: refresh Eligible/Registered/Total_Stock for the
#    depositories whose latest-day figures changed
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Today_Summary")

# BRINK'S, INC. (row 3)
$wsSummary.Range("B3").Value = 23710.274
$wsSummary.Range("C3").Value = 95517.77499999999
$wsSummary.Range("D3").Value = 119228.049

# DELAWARE DEPOSITORY (row 5)
$wsSummary.Range("B5").Value = 18459.584
$wsSummary.Range("D5").Value = 20093.525

# MANFRA, TORDELLA & BROOKES, LLC (row 11)
$wsSummary.Range("C11").Value = 60301.249
$wsSummary.Range("D11").Value = 61369.65700000001

# ---------------------------------------------------------------------
# 3) Monthly_Stats: refresh the 2026-01 grand totals and the monthly
#    detail rows for the same depositories
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")

# 2026-01 grand total (row 2)
$wsMonthly.Range("B2").Value = 294309.517
$wsMonthly.Range("C2").Value = 369983.665
$wsMonthly.Range("D2").Value = 664293.182

# BRINK'S, INC. Eligible / Registered detail (rows 9-10)
$wsMonthly.Range("E9").Value = 23710.274
$wsMonthly.Range("E10").Value = 95517.77499999999

# DELAWARE DEPOSITORY Eligible detail (row 13)
$wsMonthly.Range("D13").Value = 50.145
$wsMonthly.Range("E13").Value = 18459.584

# MANFRA, TORDELLA & BROOKES, LLC Registered detail (row 26)
$wsMonthly.Range("E26").Value = 60301.249
